$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F, matching style of existing headers (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "連続・大幅利下げ"

# New values for column F (rows 2-13)
$values = @(
    5.374999999999997,
    5.375000000000038,
    5.324995262057104,
    4.00851409442404,
    4.061259621418382,
    4.060793400635059,
    3.994474611393215,
    3.920048237798839,
    3.820634323769962,
    3.698741862372052,
    3.560296115450755,
    3.410189600472426
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
